$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values must be swapped between row 4 and row 6
$cols = @("A", "B", "D", "E", "F", "G", "H", "S", "AC", "AI")

foreach ($col in $cols) {
    $addr4 = $col + "4"
    $addr6 = $col + "6"
    $val4 = $ws.Range($addr4).Value2
    $val6 = $ws.Range($addr6).Value2
    $ws.Range($addr4).Value2 = $val6
    $ws.Range($addr6).Value2 = $val4
}
